$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shared strings "TestSignupa14".."TestSignupa17" were renamed to "TestSignupc14".."TestSignupc17".
# These values live in columns A and D, rows 2-5 (A2/D2, A3/D3, A4/D4, A5/D5).
$ws.Range("A2").Value = "TestSignupc17"
$ws.Range("D2").Value = "TestSignupc17"

$ws.Range("A3").Value = "TestSignupc14"
$ws.Range("D3").Value = "TestSignupc14"

$ws.Range("A4").Value = "TestSignupc15"
$ws.Range("D4").Value = "TestSignupc15"

$ws.Range("A5").Value = "TestSignupc16"
$ws.Range("D5").Value = "TestSignupc16"
